# edit.ps1
#
# Reproduces the content-level edits captured by the target diff:
#
#   1. Slide 2  ("MAUDE"): the acronym line "Manufacturer and User
#      Facility Device Experience" had the runs "ser " / "F" (bold) /
#      "acility " collapsed into a single, non-bold run "ser Facility ".
#   2. Slide 12 ("Approach II - Data Cleaning"): the bullet "Selected 5
#      major metadata attributes to classify data " + "on" (two runs)
#      was collapsed into a single run "...classify data on".
#   3. Slide 13 ("Approach II - Analysis and Results"): the bullet
#      "76.67" + "% positive " + "correlation" (three runs) was
#      collapsed into a single run "76.67% positive correlation".

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Slide 2 - "Manufacturer and User Facility Device Experience"
#    Replace the 13 characters "ser Facility " (spanning the old
#    "ser " / "F" / "acility " runs) with themselves so the host
#    collapses them into one plain run, matching the diff.
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$shape2 = $s2.Shapes.Item(2)
$tr2 = $shape2.TextFrame.TextRange
$para2 = $tr2.Paragraphs(1, 1)
$sub2 = $para2.Characters(19, 13)
$sub2.Text = "ser Facility "

# ---------------------------------------------------------------------
# 2) Slide 12 - "Selected 5 major metadata attributes to classify data on"
#    Round-trip the paragraph text through a placeholder so the two
#    runs get rewritten as a single run with the same final text.
# ---------------------------------------------------------------------
$s12 = $p.Slides.Item(12)
$shape12 = $s12.Shapes.Item(3)
$tr12 = $shape12.TextFrame.TextRange
$para12 = $tr12.Paragraphs(1, 1)
$para12.Text = "zzz placeholder zzz"
$para12.Text = "Selected 5 major metadata attributes to classify data on"

# ---------------------------------------------------------------------
# 3) Slide 13 - "76.67% positive correlation"
#    Same round-trip technique to collapse the three runs into one.
# ---------------------------------------------------------------------
$s13 = $p.Slides.Item(13)
$shape13 = $s13.Shapes.Item(3)
$tr13 = $shape13.TextFrame.TextRange
$para13 = $tr13.Paragraphs(3, 1)
$para13.Text = "zzz placeholder zzz"
$para13.Text = "76.67% positive correlation"
